# Update "want to attend" counts (column F) on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets to match the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 6).Value  = 37
$ws1.Cells.Item(7, 6).Value  = 4332
$ws1.Cells.Item(14, 6).Value = 175
$ws1.Cells.Item(15, 6).Value = 956
$ws1.Cells.Item(16, 6).Value = 73
$ws1.Cells.Item(19, 6).Value = 66
$ws1.Cells.Item(22, 6).Value = 3441
$ws1.Cells.Item(23, 6).Value = 5784
$ws1.Cells.Item(27, 6).Value = 520
$ws1.Cells.Item(29, 6).Value = 3343
$ws1.Cells.Item(30, 6).Value = 355
$ws1.Cells.Item(31, 6).Value = 22
$ws1.Cells.Item(32, 6).Value = 2454
$ws1.Cells.Item(34, 6).Value = 516
$ws1.Cells.Item(35, 6).Value = 121
$ws1.Cells.Item(40, 6).Value = 1007
$ws1.Cells.Item(41, 6).Value = 894
$ws1.Cells.Item(44, 6).Value = 46
$ws1.Cells.Item(45, 6).Value = 39
$ws1.Cells.Item(46, 6).Value = 464
$ws1.Cells.Item(48, 6).Value = 545

# --- Sheet 4: 全部类型 ------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2, 6).Value  = 37
$ws4.Cells.Item(7, 6).Value  = 4332
$ws4.Cells.Item(15, 6).Value = 175
$ws4.Cells.Item(16, 6).Value = 956
$ws4.Cells.Item(17, 6).Value = 73
$ws4.Cells.Item(20, 6).Value = 66
$ws4.Cells.Item(23, 6).Value = 3441
$ws4.Cells.Item(24, 6).Value = 5784
$ws4.Cells.Item(28, 6).Value = 520
$ws4.Cells.Item(30, 6).Value = 3343
$ws4.Cells.Item(31, 6).Value = 355
$ws4.Cells.Item(32, 6).Value = 22
$ws4.Cells.Item(33, 6).Value = 2454
$ws4.Cells.Item(35, 6).Value = 516
$ws4.Cells.Item(36, 6).Value = 121
$ws4.Cells.Item(41, 6).Value = 1007
$ws4.Cells.Item(42, 6).Value = 894
$ws4.Cells.Item(45, 6).Value = 46
$ws4.Cells.Item(46, 6).Value = 39
$ws4.Cells.Item(47, 6).Value = 464
$ws4.Cells.Item(49, 6).Value = 545
